$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3619.3684
$ws.Range("J17").Value = 3784.4666
$ws.Range("L17").Value = 11353.3998
$ws.Range("N17").Value = -11689.3998
$ws.Range("H96").Value = 723.75
$ws.Range("I96").Value = 633.3333
$ws.Range("J96").Value = 995
$ws.Range("K96").Value = 1899.9999
$ws.Range("L96").Value = 2985
$ws.Range("M96").Value = -526.9999
$ws.Range("N96").Value = -5731
$ws.Range("H111").Value = 553
$ws.Range("I111").Value = 538.0769
$ws.Range("K111").Value = 1614.2307
$ws.Range("M111").Value = 1452.7693
$ws.Range("H127").Value = 823.8570999999999
$ws.Range("I127").Value = 823.8570999999999
$ws.Range("K127").Value = 2471.5713
$ws.Range("M127").Value = 2488.4287
$ws.Range("H131").Value = 316.5
$ws.Range("I131").Value = 329.8
$ws.Range("K131").Value = 989.4000000000001
$ws.Range("M131").Value = 4050.6
$ws.Range("H132").Value = 4851
$ws.Range("I132").Value = 3108.3635
$ws.Range("J132").Value = 24020
$ws.Range("K132").Value = 9325.0905
$ws.Range("L132").Value = 72060
$ws.Range("M132").Value = -6795.0905
$ws.Range("N132").Value = -77120
$ws.Range("H134").Value = 80000
$ws.Range("J134").Value = 80000
$ws.Range("L134").Value = 80000
$ws.Range("N134").Value = -90140

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2174.4443
$ws.Range("I2").Value = 1581.4286
$ws.Range("J2").Value = 4250
$ws.Range("K2").Value = 1581.4286
$ws.Range("L2").Value = 4250
$ws.Range("M2").Value = -1468.4286
$ws.Range("N2").Value = -4476
$ws.Range("H25").Value = 2768.6667
$ws.Range("I25").Value = 612.5
$ws.Range("J25").Value = 20018
$ws.Range("K25").Value = 612.5
$ws.Range("L25").Value = 20018
$ws.Range("M25").Value = -210.5
$ws.Range("N25").Value = -20822
$ws.Range("H45").Value = 1675.2142
$ws.Range("I45").Value = 1688.6923
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 1688.6923
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -1311.6923
$ws.Range("N45").Value = -2254
$ws.Range("H74").Value = 2428
$ws.Range("I74").Value = 3856
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 3856
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -2982
$ws.Range("N74").Value = -2748
$ws.Range("H77").Value = 2428
$ws.Range("I77").Value = 3856
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 19280
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -14912
$ws.Range("N77").Value = -13736
$ws.Range("H88").Value = 3524.9092
$ws.Range("I88").Value = 1443.75
$ws.Range("J88").Value = 4714.143
$ws.Range("K88").Value = 1443.75
$ws.Range("L88").Value = 4714.143
$ws.Range("M88").Value = -1037.75
$ws.Range("N88").Value = -5526.143
$ws.Range("H91").Value = 3524.9092
$ws.Range("I91").Value = 1443.75
$ws.Range("J91").Value = 4714.143
$ws.Range("K91").Value = 1443.75
$ws.Range("L91").Value = 4714.143
$ws.Range("M91").Value = -39.75
$ws.Range("N91").Value = -7522.143
$ws.Range("H116").Value = 2174.4443
$ws.Range("I116").Value = 1581.4286
$ws.Range("J116").Value = 4250
$ws.Range("K116").Value = 1581.4286
$ws.Range("L116").Value = 4250
$ws.Range("M116").Value = 712.5714
$ws.Range("N116").Value = -8838

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2174.4443
$ws.Range("I3").Value = 1581.4286
$ws.Range("J3").Value = 4250
$ws.Range("K3").Value = 1581.4286
$ws.Range("L3").Value = 4250
$ws.Range("M3").Value = -1467.4286
$ws.Range("N3").Value = -4478
$ws.Range("H16").Value = 198.875
$ws.Range("I16").Value = 200
$ws.Range("J16").Value = 198.71428
$ws.Range("K16").Value = 200
$ws.Range("L16").Value = 198.71428
$ws.Range("M16").Value = -30
$ws.Range("N16").Value = -538.71428
$ws.Range("H94").Value = 2440.2727
$ws.Range("I94").Value = 2554.3
$ws.Range("K94").Value = 2554.3
$ws.Range("M94").Value = -2103.3
$ws.Range("H134").Value = 3000.0645
$ws.Range("I134").Value = 3000.0645
$ws.Range("K134").Value = 9000.193499999999
$ws.Range("M134").Value = -6465.193499999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6499.5
$ws.Range("I31").Value = 6499.5
$ws.Range("K31").Value = 6499.5
$ws.Range("M31").Value = -6204.5
$ws.Range("H34").Value = 6499.5
$ws.Range("I34").Value = 6499.5
$ws.Range("K34").Value = 6499.5
$ws.Range("M34").Value = -6297.5
$ws.Range("H41").Value = 19297.5
$ws.Range("J41").Value = 19996.875
$ws.Range("L41").Value = 19996.875
$ws.Range("N41").Value = -20852.875
$ws.Range("H50").Value = 29996.875
$ws.Range("J50").Value = 29996.875
$ws.Range("L50").Value = 29996.875
$ws.Range("N50").Value = -31246.875
$ws.Range("H60").Value = 19428.572
$ws.Range("H132").Value = 4666.6665
$ws.Range("I132").Value = 5250
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 15750
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -13220
$ws.Range("N132").Value = -15560

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 178.81818
$ws.Range("J12").Value = 256.57144
$ws.Range("L12").Value = 769.71432
$ws.Range("N12").Value = -1115.71432
$ws.Range("H123").Value = 2998.1428
$ws.Range("I123").Value = 2995.6667
$ws.Range("J123").Value = 3000
$ws.Range("K123").Value = 8987.000100000001
$ws.Range("L123").Value = 9000
$ws.Range("M123").Value = -6537.000100000001
$ws.Range("N123").Value = -13900
$ws.Range("H124").Value = 1975
$ws.Range("I124").Value = 1975
$ws.Range("K124").Value = 5925
$ws.Range("M124").Value = -1015
$ws.Range("H125").Value = 5000
$ws.Range("I125").Value = 5000
$ws.Range("K125").Value = 15000
$ws.Range("M125").Value = -10080
$ws.Range("H126").Value = 588
$ws.Range("I126").Value = 560
$ws.Range("K126").Value = 1680
$ws.Range("M126").Value = 3260
$ws.Range("H131").Value = 1775.7778
$ws.Range("I131").Value = 999.5
$ws.Range("J131").Value = 1997.5714
$ws.Range("K131").Value = 2998.5
$ws.Range("L131").Value = 5992.7142
$ws.Range("M131").Value = 2041.5
$ws.Range("N131").Value = -16072.7142

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 5000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H132").Value = 3512.8333
$ws.Range("I132").Value = 2026
$ws.Range("K132").Value = 6078
$ws.Range("M132").Value = -3548

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4000
$ws.Range("I7").Value = 4000
$ws.Range("K7").Value = 4000
$ws.Range("M7").Value = -3888
$ws.Range("H16").Value = 1478.2
$ws.Range("I16").Value = 1478.2
$ws.Range("K16").Value = 1478.2
$ws.Range("M16").Value = -1308.2
$ws.Range("H122").Value = 3493
$ws.Range("I122").Value = 3493
$ws.Range("K122").Value = 10479
$ws.Range("M122").Value = -8029
$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530
$ws.Range("H135").Value = 42000
$ws.Range("J135").Value = 42000
$ws.Range("L135").Value = 42000
$ws.Range("N135").Value = -52140
$ws.Range("H136").Value = 3094.9167
$ws.Range("I136").Value = 3143.4
$ws.Range("J136").Value = 2852.5
$ws.Range("K136").Value = 9430.200000000001
$ws.Range("L136").Value = 8557.5
$ws.Range("M136").Value = -6880.200000000001
$ws.Range("N136").Value = -13657.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("N135").Value = -70140
